$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to Text format so numeric-looking strings
# (e.g. "1.001", "237.76") are stored as text, matching the source data,
# not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

# Price (D) updates
$ws.Range("D2").Value = "26.086.47"
$ws.Range("D3").Value = "1.767.44"
$ws.Range("D5").Value = "237.76"
$ws.Range("D6").Value = "1.001"
$ws.Range("D7").Value = "0.5236"
$ws.Range("D8").Value = "0.2756"
$ws.Range("D9").Value = "40.41"
$ws.Range("D10").Value = "0.06196"
$ws.Range("D11").Value = "1.777.72"
$ws.Range("D12").Value = "0.07018"
$ws.Range("D13").Value = "15.77"
$ws.Range("D14").Value = "0.6448"
$ws.Range("D15").Value = "4.541"
$ws.Range("D16").Value = "78.11"
$ws.Range("D17").Value = "1.001"
$ws.Range("D18").Value = "1.001"
$ws.Range("D19").Value = "26.109.42"
$ws.Range("D20").Value = "11.65"
$ws.Range("D21").Value = "0.000006744"
$ws.Range("D22").Value = "2.002.80"
$ws.Range("D24").Value = "8.455"
$ws.Range("D25").Value = "5.189"
$ws.Range("D26").Value = "138.79"
$ws.Range("D27").Value = "1.487"
$ws.Range("D28").Value = "1.856"
$ws.Range("D29").Value = "15.17"
$ws.Range("D30").Value = "103.29"
$ws.Range("D31").Value = "0.08408"
$ws.Range("D32").Value = "3.698"
$ws.Range("D33").Value = "3.453"
$ws.Range("D34").Value = "0.04453"
$ws.Range("D35").Value = "2.594"
$ws.Range("D36").Value = "1.005"
$ws.Range("D37").Value = "0.6052"
$ws.Range("D39").Value = "0.01591"
$ws.Range("D40").Value = "1.986"
$ws.Range("D42").Value = "102.71"
$ws.Range("D43").Value = "0.3880"
$ws.Range("D44").Value = "0.7504"
$ws.Range("D45").Value = "4.937"
$ws.Range("D47").Value = "6.365"
$ws.Range("D48").Value = "0.1119"
$ws.Range("D49").Value = "30.20"
$ws.Range("D50").Value = "52.64"

# Volume(1h) (E) updates
$ws.Range("E2").Value = "  +1.19%  "
$ws.Range("E3").Value = "  +1.24%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("E5").Value = "  -0.27%  "
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("E7").Value = "  +3.51%  "
$ws.Range("E8").Value = "  +0.95%  "
$ws.Range("E9").Value = "  -3.81%  "
$ws.Range("E10").Value = "  +0.69%  "
$ws.Range("E11").Value = "  +1.85%  "
$ws.Range("E12").Value = "  +1.27%  "
$ws.Range("E13").Value = "  +1.62%  "
$ws.Range("E14").Value = "  +7.21%  "
$ws.Range("E15").Value = "  +0.31%  "
$ws.Range("E17").Value = "  +0.06%  "
$ws.Range("E18").Value = "  +0.13%  "
$ws.Range("E19").Value = "  +1.26%  "
$ws.Range("E20").Value = "  -0.02%  "
$ws.Range("E21").Value = "  -2.17%  "
$ws.Range("E22").Value = "  +1.73%  "
$ws.Range("E23").Value = "  +0.22%  "
$ws.Range("E24").Value = "  +3.46%  "
$ws.Range("E25").Value = "  -1.05%  "
$ws.Range("E26").Value = "  +0.62%  "
$ws.Range("E27").Value = "  +1.35%  "
$ws.Range("E28").Value = "  +1.94%  "
$ws.Range("E29").Value = "  +1.13%  "
$ws.Range("E30").Value = "  -0.61%  "
$ws.Range("E31").Value = "  +3.35%  "
$ws.Range("E32").Value = "  -0.44%  "
$ws.Range("E33").Value = "  -1.04%  "
$ws.Range("E34").Value = "  -2.01%  "
$ws.Range("E35").Value = "  -0.81%  "
$ws.Range("E36").Value = "  +1.97%  "
$ws.Range("E37").Value = "  -0.86%  "
$ws.Range("E38").Value = "  +3.08%  "
$ws.Range("E39").Value = "  +2.28%  "
$ws.Range("E40").Value = "  +2.79%  "
$ws.Range("E41").Value = "  +0.27%  "
$ws.Range("E42").Value = "  +0.72%  "
$ws.Range("E43").Value = "  +0.75%  "
$ws.Range("E45").Value = "  -0.54%  "
$ws.Range("E46").Value = "  +2.51%  "
$ws.Range("E47").Value = "  +6.54%  "
$ws.Range("E48").Value = "  +0.55%  "
$ws.Range("E49").Value = "  +0.08%  "
$ws.Range("E50").Value = "  +0.14%  "
$ws.Range("E51").Value = "  +0.71%  "

# Restore default style on column D so only the values (and text-type) changed,
# without leaving a lingering custom number format on the cells.
$ws.Range("D2:D51").Style = "Normal"
